$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 36 and 37, pushing the existing rows 36-52 down to 38-54.
$ws.Rows("36:37").Insert()

# Column D values look numeric ("0.086", "0.916", "0.120") - force text
# formatting first so Excel's autodetection doesn't coerce them to numbers
# and strip the padding / trailing zero, matching the original sheet where
# these are stored as literal padded strings.
$ws.Range("D35:D37").NumberFormat = "@"

# Row 35: Pseudotime -> Pseudotime_1 (new stats)
$ws.Range("A35").Value = "Pseudotime_1"
$ws.Range("B35").Value = "  5.67 (3.78)  "
$ws.Range("C35").Value = "  4.63 (3.33)  "
$ws.Range("D35").Value = "  0.086  "

# Row 36: new Pseudotime_2 row
$ws.Range("A36").Value = "Pseudotime_2"
$ws.Range("B36").Value = "  4.63 (2.81)  "
$ws.Range("C36").Value = "  4.58 (3.13)  "
$ws.Range("D36").Value = "  0.916  "

# Row 37: new Pseudotime_3 row
$ws.Range("A37").Value = "Pseudotime_3"
$ws.Range("B37").Value = "  4.35 (2.84)  "
$ws.Range("C37").Value = "  3.67 (2.26)  "
$ws.Range("D37").Value = "  0.120  "
